$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Format" column (G) was left blank for a handful of rows whose
# transition/duration bookkeeping depends on it. Backfill those rows with
# the generic "Событие" ("Event") format so downstream transition-time
# calculations for the route list have a value to read.
$ws.Range("G3").Value = "Событие"
$ws.Range("G13").Value = "Событие"
$ws.Range("G21").Value = "Событие"
$ws.Range("G22").Value = "Событие"

# Leave the selection where the author ended up after the fix.
$ws.Range("G14").Select()
